$wb = $excel.ActiveWorkbook

# --- "2-data" sheet: fix the tag-group id (A1: 2 -> 9) and rename the sheet to match ---
$wsData = $wb.Worksheets.Item("2-data")
$wsData.Range("A1").Value = 9
$wsData.Name = "9-data"

# --- Make "9-data" the active sheet, with B16:B17 selected on it ---
$wsData.Select()
$wsData.Range("B16:B17").Select()
